# Add a new "Skill Description" column after column A (SkillCode),
# shifting SFIA Level / Keycode / Description one column to the right,
# and populate the new column with the skill's full/display name
# (falling back to the skill code itself when no full name is known yet).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank column at B; existing B/C/D move to C/D/E.
$ws.Columns("B:B").Insert()

# Header
$ws.Range("B1").Value = "Skill Description"

# Map of SkillCode (column A) -> full Skill Description (new column B).
# Rows with a known full name get it; everything else repeats the code.
$descriptions = @{
    "USEV" = "User experience evaluation"
}

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
for ($r = 2; $r -le $lastRow; $r++) {
    $code = $ws.Cells.Item($r, 1).Value2
    if ($descriptions.ContainsKey($code)) {
        $ws.Cells.Item($r, 2).Value = $descriptions[$code]
    } else {
        $ws.Cells.Item($r, 2).Value = $code
    }
}
